# Open registration for NextGen Event: update Ottavia Prunas' company and
# add a new roster entry for Nikki Rommers.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2020 Roster")

# New roster row for Nikki Rommers at University Hospital of Basel,
# opted into the mailing list / committee (column E).
$ws.Range("A11").Value = "Nikki"
$ws.Range("B11").Value = "Rommers"
$ws.Range("C11").Value = "University Hospital of Basel"
$ws.Range("E11").Value = 1

# Ottavia Prunas moved from "Swiss Tropical and Public Health Institute"
# to "University Hospital of Basel"
$ws.Range("C9").Value = "University Hospital of Basel"

# Match the styling used by the other data rows (row height + cell styles).
$ws.Range("A10:C10").Copy()
$ws.Range("A11:C11").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("E10").Copy()
$ws.Range("E11").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = $false
$ws.Rows.Item(11).RowHeight = 16

# Update the selected cell shown when the sheet was last saved.
$ws.Range("C11").Select()
